# Swap the two product names in column B (row 2 <-> row 4) and update
# the active selection to B2:B5 with B2 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B2").Value = "Pril ISIS Cold Power liquid 650ml Lemon"
$ws.Range("B4").Value = "Bref 900ml javel disinf"

$ws.Range("B2:B5").Select()
